$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026888833949624
$ws.Range("D2").Value = 1.034688655946679
$ws.Range("E2").Value = 1.027055598426267
$ws.Range("F2").Value = 1.041745439449038
$ws.Range("I2").Value = 1.032308198255865
$ws.Range("J2").Value = 1.032049764430267
$ws.Range("K2").Value = 1.03748741061062
$ws.Range("L2").Value = 1.029876442648579
$ws.Range("M2").Value = 1.044524084905264
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027866144984695
$ws.Range("D3").Value = 1.035436508606351
$ws.Range("E3").Value = 1.027885583570498
$ws.Range("F3").Value = 1.042687888619449
$ws.Range("I3").Value = 1.03248595447643
$ws.Range("J3").Value = 1.032666883976672
$ws.Range("K3").Value = 1.038044666466241
$ws.Range("L3").Value = 1.030514011875562
$ws.Range("M3").Value = 1.045276878749218
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028498742283245
$ws.Range("D4").Value = 1.035920201214203
$ws.Range("E4").Value = 1.028423216254621
$ws.Range("F4").Value = 1.043298023170621
$ws.Range("I4").Value = 1.032599109429178
$ws.Range("J4").Value = 1.03306582938288
$ws.Range("K4").Value = 1.038404367731821
$ws.Range("L4").Value = 1.030926494089155
$ws.Range("M4").Value = 1.04576367094316
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028764736190355
$ws.Range("D5").Value = 1.036123492101626
$ws.Range("E5").Value = 1.028649373964757
$ws.Range("F5").Value = 1.043554596009449
$ws.Range("I5").Value = 1.032646232566291
$ws.Range("J5").Value = 1.033233456081535
$ws.Range("K5").Value = 1.038555374351527
$ws.Range("L5").Value = 1.031099884506132
$ws.Range("M5").Value = 1.045968241702309
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.0288094006502
$ws.Range("D6").Value = 1.036157622362469
$ws.Range("E6").Value = 1.028687354835367
$ws.Range("F6").Value = 1.043597679921151
$ws.Range("I6").Value = 1.03265411850446
$ws.Range("J6").Value = 1.033261596042534
$ws.Range("K6").Value = 1.038580716578932
$ws.Range("L6").Value = 1.031128996490909
$ws.Range("M6").Value = 1.046002585504774
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028502296311227
$ws.Range("D7").Value = 1.035922917808771
$ws.Range("E7").Value = 1.028426237648432
$ws.Range("F7").Value = 1.043301451226987
$ws.Range("I7").Value = 1.032599740848904
$ws.Range("J7").Value = 1.033068069573541
$ws.Range("K7").Value = 1.03840638632312
$ws.Range("L7").Value = 1.030928811008739
$ws.Range("M7").Value = 1.045766404730327
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02721907654529
$ws.Range("D8").Value = 1.034941440640783
$ws.Range("E8").Value = 1.027335975668092
$ws.Range("F8").Value = 1.042063880287147
$ws.Range("I8").Value = 1.032368657604953
$ws.Range("J8").Value = 1.032258399432114
$ws.Range("K8").Value = 1.037675919563111
$ws.Range("L8").Value = 1.030091925635986
$ws.Range("M8").Value = 1.044778559838645
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024959530446505
$ws.Range("D9").Value = 1.033210338294153
$ws.Range("E9").Value = 1.025419263077521
$ws.Range("F9").Value = 1.039885525481526
$ws.Range("I9").Value = 1.031947203399675
$ws.Range("J9").Value = 1.030828843421801
$ws.Range("K9").Value = 1.036382043242561
$ws.Range("L9").Value = 1.028616747355215
$ws.Range("M9").Value = 1.043035480055882
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023454314194847
$ws.Range("D10").Value = 1.032055258775724
$ws.Range("E10").Value = 1.024144526024411
$ws.Range("F10").Value = 1.038434974015943
$ws.Range("I10").Value = 1.031656687525171
$ws.Range("J10").Value = 1.029873962585304
$ws.Range("K10").Value = 1.035515010905864
$ws.Range("L10").Value = 1.02763302454747
$ws.Range("M10").Value = 1.041871895214608
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022802818411989
$ws.Range("D11").Value = 1.0315548718669
$ws.Range("E11").Value = 1.023593293024358
$ws.Range("F11").Value = 1.037807282919505
$ws.Range("I11").Value = 1.031528634482487
$ws.Range("J11").Value = 1.029460061131233
$ws.Range("K11").Value = 1.035138532855439
$ws.Range("L11").Value = 1.027207008016315
$ws.Range("M11").Value = 1.041367699481475
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022560865349373
$ws.Range("D12").Value = 1.031368972522702
$ws.Range("E12").Value = 1.023388652268683
$ws.Range("F12").Value = 1.037574192752014
$ws.Range("I12").Value = 1.031480731256633
$ws.Range("J12").Value = 1.029306255566066
$ws.Range("K12").Value = 1.034998535667919
$ws.Range("L12").Value = 1.027048758583374
$ws.Range("M12").Value = 1.041180366235033
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022612763217662
$ws.Range("D13").Value = 1.031408850054329
$ws.Range("E13").Value = 1.023432543332632
$ws.Range("F13").Value = 1.037624188549985
$ws.Range("I13").Value = 1.031491021975425
$ws.Range("J13").Value = 1.029339250272519
$ws.Range("K13").Value = 1.035028572606837
$ws.Range("L13").Value = 1.027082703969761
$ws.Range("M13").Value = 1.041220552212835
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022782817647329
$ws.Range("D14").Value = 1.031539506043062
$ws.Range("E14").Value = 1.023576375070371
$ws.Range("F14").Value = 1.037788014337435
$ws.Range("I14").Value = 1.031524681691402
$ws.Range("J14").Value = 1.029447348819388
$ws.Range("K14").Value = 1.035126963831304
$ws.Range("L14").Value = 1.027193927225881
$ws.Range("M14").Value = 1.041352215522396
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022887599305153
$ws.Range("D15").Value = 1.031620003123302
$ws.Range("E15").Value = 1.023665009387368
$ws.Range("F15").Value = 1.037888961093081
$ws.Range("I15").Value = 1.031545375699527
$ws.Range("J15").Value = 1.029513943424917
$ws.Range("K15").Value = 1.035187565203304
$ws.Range("L15").Value = 1.027262454518513
$ws.Range("M15").Value = 1.04143333070436
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023497557613938
$ws.Range("D16").Value = 1.03208846306383
$ws.Range("E16").Value = 1.024181125152565
$ws.Range("F16").Value = 1.038476640478487
$ws.Range("I16").Value = 1.031665138465034
$ws.Range("J16").Value = 1.029901422827514
$ws.Range("K16").Value = 1.035539974522761
$ws.Range("L16").Value = 1.027661296706492
$ws.Range("M16").Value = 1.041905349609902
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023880241729207
$ws.Range("D17").Value = 1.032382254895039
$ws.Range("E17").Value = 1.024505068886561
$ws.Range("F17").Value = 1.038845386020944
$ws.Range("I17").Value = 1.031739658428553
$ws.Range("J17").Value = 1.030144363343931
$ws.Range("K17").Value = 1.035760751796629
$ws.Range("L17").Value = 1.027911464893266
$ws.Range("M17").Value = 1.042201340012505
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024103481185913
$ws.Range("D18").Value = 1.032553596444101
$ws.Range("E18").Value = 1.024694090781683
$ws.Range("F18").Value = 1.039060508280587
$ws.Range("I18").Value = 1.031782906783699
$ws.Range("J18").Value = 1.030286024730531
$ws.Range("K18").Value = 1.035889426335164
$ws.Range("L18").Value = 1.028057378108424
$ws.Range("M18").Value = 1.042373951811446
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024179604515228
$ws.Range("D19").Value = 1.032612015673954
$ws.Range("E19").Value = 1.024758554372811
$ws.Range("F19").Value = 1.039133866032292
$ws.Range("I19").Value = 1.031797616381509
$ws.Range("J19").Value = 1.03033432050613
$ws.Range("K19").Value = 1.035933283850922
$ws.Range("L19").Value = 1.028107129760413
$ws.Range("M19").Value = 1.042432802111033
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023839180622688
$ws.Range("D20").Value = 1.03235073610547
$ws.Range("E20").Value = 1.024470305438503
$ws.Range("F20").Value = 1.038805819049164
$ws.Range("I20").Value = 1.031731685673731
$ws.Range("J20").Value = 1.030118302434245
$ws.Range("K20").Value = 1.035737074936483
$ws.Range("L20").Value = 1.027884624793565
$ws.Range("M20").Value = 1.042169586594206
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022732739707282
$ws.Range("D21").Value = 1.031501032022638
$ws.Range("E21").Value = 1.02353401711585
$ws.Range("F21").Value = 1.037739769999224
$ws.Range("I21").Value = 1.031514779087063
$ws.Range("J21").Value = 1.02941551825625
$ws.Range("K21").Value = 1.035097994380585
$ws.Range("L21").Value = 1.027161174961488
$ws.Range("M21").Value = 1.041313445401006
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022037316619427
$ws.Range("D22").Value = 1.030966596865333
$ws.Range("E22").Value = 1.022945982412403
$ws.Range("F22").Value = 1.037069863355898
$ws.Range("I22").Value = 1.031376442414379
$ws.Range("J22").Value = 1.028973279236845
$ws.Range("K22").Value = 1.034695274170843
$ws.Range("L22").Value = 1.026706267880496
$ws.Range("M22").Value = 1.040774852163037
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022405950619563
$ws.Range("D23").Value = 1.031249928910784
$ws.Range("E23").Value = 1.023257649058972
$ws.Range("F23").Value = 1.037424958996634
$ws.Range("I23").Value = 1.031449962790701
$ws.Range("J23").Value = 1.029207753405105
$ws.Range("K23").Value = 1.034908849257943
$ws.Range("L23").Value = 1.026947426774364
$ws.Range("M23").Value = 1.041060399059458
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023857734288
$ws.Range("D24").Value = 1.032364978158609
$ws.Range("E24").Value = 1.024486013325014
$ws.Range("F24").Value = 1.038823697539425
$ws.Range("I24").Value = 1.031735288891686
$ws.Range("J24").Value = 1.030130078366441
$ws.Range("K24").Value = 1.035747773802876
$ws.Range("L24").Value = 1.027896752696816
$ws.Range("M24").Value = 1.04218393470344
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025543477721435
$ws.Range("D25").Value = 1.03365805338119
$ws.Range("E25").Value = 1.025914242776906
$ws.Range("F25").Value = 1.04044839045619
$ws.Range("I25").Value = 1.032057844835145
$ws.Range("J25").Value = 1.031198746368658
$ws.Range("K25").Value = 1.036717329444483
$ws.Range("L25").Value = 1.028998167766765
$ws.Range("M25").Value = 1.043486382472114
